# Remove the trailing "Ver no Jupiter..." / "(c) 2020 ..." footer block
# (and the blank paragraph that precedes it) that followed the last
# "Requisitos" entry ("LOQ4086: Operacoes Unitarias II (Requisito fraco)").

$d = $word.ActiveDocument

# Locate the anchor paragraphs by their text so the script is resilient to
# any incidental shift in paragraph indices.
$startIndex = $null
$endIndex = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*LOQ4086*") {
        $startIndex = $i + 1
    }
    if (($t -like "*Ver no Jupiter*") -or ($t -like "*Creative Commons Attribution*")) {
        $endIndex = $i
    }
}

if (($startIndex -ne $null) -and ($endIndex -ne $null) -and ($startIndex -le $endIndex)) {
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($endIndex)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()

    Write-Output "Removed paragraphs $startIndex through $endIndex."
} else {
    Write-Output "Anchor paragraphs not found; no changes made."
}
